$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "boolean" column header in I1 (added to shared strings first, matching
# the original author's string table order)
$ws.Range("I1").Value = "boolean[boolean]"
$ws.Range("I1").WrapText = $true

# Boolean values as literal text "true"/"false" (not real TRUE/FALSE booleans).
# A direct .Value assignment of "true"/"false" auto-converts to a boolean
# cell, so instead build the text via a formula and paste the computed
# value back in as plain text.
$ws.Range("K1").Formula = "=""true"""
$ws.Range("K1").Copy()
$ws.Range("I2").PasteSpecial(-4163)

$ws.Range("K1").Formula = "=""false"""
$ws.Range("K1").Copy()
$ws.Range("I3").PasteSpecial(-4163)

$ws.Range("K1").ClearContents()

# New column width for the boolean column
$ws.Columns.Item(9).ColumnWidth = 19.5

# Row 1 grows taller to fit the wrapped header text
$ws.Rows.Item(1).RowHeight = 17

# Selection ends up on I4 after the edits
$ws.Range("I4").Select()
